# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E35) listed the 20 overdue periods in
# descending order (2108 ... 2001). The updated account-statement database
# lists them in ascending order (2001 ... 2108) instead, so every row's
# period value is reversed top-to-bottom. The two rows whose "Valor Mora"
# (F column) differed from the common 80000 value - originally 77334 on the
# 2108 row and 74667 on the 2001 row - travel along with their period, so
# after the reorder F16 (now period 2001) holds 74667 and F35 (now period
# 2108) holds 77334.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 16-35.
$periods = @(2001, 2002, 2003, 2004, 2005, 2006, 2007, 2008, 2009, 2010, `
             2011, 2012, 2101, 2102, 2103, 2104, 2105, 2106, 2107, 2108)

# Valor Mora amounts matched to the same ascending period order (only the
# first and last rows differ from 80000, mirroring the original data).
$valores = @(74667, 80000, 80000, 80000, 80000, 80000, 80000, 80000, 80000, 80000, `
             80000, 80000, 80000, 80000, 80000, 80000, 80000, 80000, 80000, 77334)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = [string]$periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
